{"js": "// Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics\n// (percentages, dollar amounts, large numbers) inside specific bullet\n// paragraphs of the resume, matching the target diff exactly.\n//\n// Strategy: for each target paragraph (identified by a unique substring),\n// locate it, then for each metric token that should be highlighted inside\n// that paragraph, use paragraph-scoped search() (matchCase + wholeWords\n// off, since tokens like \"23%\" aren't \"words\") to find the exact\n// occurrence(s) and set font.bold + font.color on the resulting range(s).\n// This naturally splits the single run into multiple runs exactly like\n// Word does when you select text and apply character formatting.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Map: a unique substring used to locate the target paragraph -> ordered\n// list of metric tokens to bold+color within that paragraph (in the order\n// they appear; each token highlighted once, first remaining match).\nconst EDITS = [\n  {\n    anchor: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    tokens: [\"23%\", \"64%\"],\n  },\n  {\n    anchor: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    tokens: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    anchor: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    tokens: [\"1,200\"],\n  },\n  {\n    anchor: \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    tokens: [\"$400M\", \"$1B\"],\n  },\n  {\n    anchor: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    tokens: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    anchor: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    tokens: [\"87%\", \"71%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Track which paragraph indices have already been consumed by an earlier\n// edit so that the duplicate \"Achieved 87% ... 71%\" short-form anchor\n// (edit #6) does not re-match the long-form paragraph (edit #2) that\n// happens to start with the same text.\nconst usedParagraphIndices = new Set();\n\nfor (const edit of EDITS) {\n  let targetParagraph = null;\n  let targetIndex = -1;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (usedParagraphIndices.has(i)) continue;\n    if (paragraphs.items[i].text === edit.anchor) {\n      targetParagraph = paragraphs.items[i];\n      targetIndex = i;\n      break;\n    }\n  }\n  if (!targetParagraph) {\n    throw new Error(\"Could not find target paragraph for anchor: \" + edit.anchor);\n  }\n  usedParagraphIndices.add(targetIndex);\n\n  for (const token of edit.tokens) {\n    const results = targetParagraph.search(token, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length === 0) {\n      throw new Error(\"Token not found in paragraph: \" + token);\n    }\n    // Highlight only the first not-yet-bolded occurrence (tokens are\n    // unique enough within each paragraph given our per-token call order).\n    const range = results.items[0];\n    range.font.bold = true;\n    range.font.color = HIGHLIGHT_COLOR;\n    await context.sync();\n  }\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics\n# (percentages, dollar amounts, large numbers) inside specific bullet\n# paragraphs of the resume, matching the target diff exactly.\n#\n# Strategy: for each target paragraph (identified by its exact full text),\n# scan $d.Paragraphs to find it, then duplicate that paragraph's Range and,\n# for each metric token (in left-to-right order), run Range.Find.Execute\n# scoped to the duplicated range. Find.Execute collapses the range to the\n# matched text and stays within the original paragraph bounds, so this is\n# equivalent to selecting just that token and applying character\n# formatting -- Word splits the run automatically around it.\n\n$d = $word.ActiveDocument\n\n# OLE (BGR) integer for RGB 2C3E50, used by Range.Font.Color.\n$HighlightColor = 5258796  # 0x503E2C == B=0x50 G=0x3E R=0x2C\n\n$Bullet = [char]0x2022\n$PlusMinus = [char]0x00B1\n$CR = [char]13\n\n# Paragraphs already matched to an earlier edit, so the short-form\n# \"Achieved 87% ... 71%\" anchor (edit #6) does not re-match the long-form\n# paragraph (edit #2) that starts with the same text.\n$usedParagraphs = @{}\n\nfunction Set-MetricHighlight($ExpectedText, $Tokens) {\n    $count = $d.Paragraphs.Count\n    $para = $null\n    for ($i = 1; $i -le $count; $i++) {\n        if ($usedParagraphs.ContainsKey($i)) {\n            continue\n        }\n        $candidate = $d.Paragraphs.Item($i)\n        $candidateText = $candidate.Range.Text.TrimEnd($CR)\n        if ($candidateText -eq $ExpectedText) {\n            $para = $candidate\n            $usedParagraphs[$i] = $true\n            break\n        }\n    }\n\n    if ($null -eq $para) {\n        throw \"Could not find target paragraph with text: [$ExpectedText]\"\n    }\n\n    foreach ($token in $Tokens) {\n        $searchRange = $para.Range.Duplicate\n        $found = $searchRange.Find.Execute($token, $true)\n        if (-not $found) {\n            throw \"Token '$token' not found in paragraph: [$ExpectedText]\"\n        }\n        $searchRange.Font.Bold = 1\n        $searchRange.Font.Color = $HighlightColor\n    }\n}\n\n$text1 = \"${Bullet} Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\nSet-MetricHighlight $text1 @(\"23%\", \"64%\")\n\n$text2 = \"${Bullet} Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ${PlusMinus}4.2% to ${PlusMinus}2.1%\"\nSet-MetricHighlight $text2 @(\"87%\", \"71%\", \"${PlusMinus}4.2%\", \"${PlusMinus}2.1%\")\n\n$text3 = \"${Bullet} Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\nSet-MetricHighlight $text3 @(\"1,200\")\n\n$text4 = \"${Bullet} Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\nSet-MetricHighlight $text4 @(\"`$400M\", \"`$1B\")\n\n$text5 = \"${Bullet} Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\nSet-MetricHighlight $text5 @(\"73.5%\", \"`$4.7M\")\n\n$text6 = \"${Bullet} Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\nSet-MetricHighlight $text6 @(\"87%\", \"71%\")\n"}
